# Commit: "Update the getmsgrecived method"
#
# massMsgSend_Data gains a duplicated creator-credential row (row 3 is
# filled in with the existing creator email/password and row 4 is a new
# clone of it), and users_LoginData gains two duplicated fan-credential
# rows (rows 3 & 4 clone row 2). The previously-active sheet
# (users_LoginData) is no longer the active tab; massMsgSend_Data becomes
# active instead, with refreshed selections on both sheets.

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("massMsgSend_Data")
$ws4 = $wb.Worksheets.Item("users_LoginData")

# --- massMsgSend_Data: row 3 already has the right style reserved for it
#     (A3/B3), just needs the creator credentials filled in; row 4 is a
#     brand-new row cloned from row 3 ---
$ws3.Range("A3").Value = "rohankapse520@gmail.com"
$ws3.Range("B3").Value = "Rohan@2025"

$ws3.Range("A3:B3").Copy() | Out-Null
$ws3.Range("A4:B4").PasteSpecial(-4122) | Out-Null
$ws3.Application.CutCopyMode = $false

$ws3.Range("A4").Value = "rohankapse520@gmail.com"
$ws3.Range("B4").Value = "Rohan@2025"

$ws3.Hyperlinks.Add($ws3.Range("A3"), "mailto:rohankapse520@gmail.com", "", "mailto:rohankapse520@gmail.com", "rohankapse520@gmail.com") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("A4"), "mailto:rohankapse520@gmail.com", "", "mailto:rohankapse520@gmail.com", "rohankapse520@gmail.com") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B3"), "mailto:Rohan@2025", "", "", "Rohan@2025") | Out-Null
$ws3.Hyperlinks.Add($ws3.Range("B4"), "mailto:Rohan@2025", "", "", "Rohan@2025") | Out-Null

# adding the hyperlinks force-applies Excel's built-in "Hyperlink" look;
# restore the original column styles (A=plain bordered, B=visited-link
# bordered) that rows 2 & 3 already use.
$ws3.Range("A2:B2").Copy() | Out-Null
$ws3.Range("A3:B3").PasteSpecial(-4122) | Out-Null
$ws3.Range("A4:B4").PasteSpecial(-4122) | Out-Null
$ws3.Application.CutCopyMode = $false

# --- users_LoginData: clone row 2's format/values into new rows 3 & 4 ---
$ws4.Range("A2:B2").Copy() | Out-Null
$ws4.Range("A3:B4").PasteSpecial(-4122) | Out-Null
$ws4.Application.CutCopyMode = $false

$ws4.Range("A3").Value = "rohan.kapse@iiclab.com"
$ws4.Range("B3").Value = "Rohan@001"
$ws4.Range("A4").Value = "rohan.kapse@iiclab.com"
$ws4.Range("B4").Value = "Rohan@001"

$ws4.Hyperlinks.Add($ws4.Range("A3"), "mailto:rohan.kapse@iiclab.com", "", "mailto:rohan.kapse@iiclab.com", "rohan.kapse@iiclab.com") | Out-Null
$ws4.Hyperlinks.Add($ws4.Range("A4"), "mailto:rohan.kapse@iiclab.com", "", "mailto:rohan.kapse@iiclab.com", "rohan.kapse@iiclab.com") | Out-Null
$ws4.Hyperlinks.Add($ws4.Range("B3"), "mailto:Rohan@001", "", "mailto:Rohan@001", "Rohan@001") | Out-Null
$ws4.Hyperlinks.Add($ws4.Range("B4"), "mailto:Rohan@001", "", "mailto:Rohan@001", "Rohan@001") | Out-Null

$ws4.Range("A2:B2").Copy() | Out-Null
$ws4.Range("A3:B3").PasteSpecial(-4122) | Out-Null
$ws4.Range("A4:B4").PasteSpecial(-4122) | Out-Null
$ws4.Application.CutCopyMode = $false

# --- switch the active sheet from users_LoginData to massMsgSend_Data,
#     and refresh each sheet's own selection ---
$ws4.Activate()
$ws4.Range("A7").Select() | Out-Null

$ws3.Activate()
$ws3.Range("A11").Select() | Out-Null
